$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -----------------
# This shared string is referenced by Overview!E2, Overview!F2, zh-cn!C2 and
# de-de!C2 - update every occurrence so they keep pointing at one shared
# string entry.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes --------------------------------------------------
# Overview: columns E (zh-cn) and F (de-de) get narrower.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C ("Status") gets narrower to match.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
